$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-01-19 Friday" "2024-01-20 Saturday"

Replace-Text "549×7=3843" "111×2=222"
Replace-Text "539×6=3234" "732×4=2928"
Replace-Text "834×8=6672" "220×9=1980"
Replace-Text "174×4=696" "760×7=5320"
Replace-Text "490×3=1470" "461×7=3227"
Replace-Text "295×5=1475" "462×8=3696"
Replace-Text "534×3=1602" "864×4=3456"
Replace-Text "260×7=1820" "614×8=4912"
Replace-Text "383×8=3064" "576×7=4032"
Replace-Text "346×3=1038" "327×6=1962"
Replace-Text "934×3=2802" "499×2=998"
Replace-Text "472×3=1416" "213×6=1278"
Replace-Text "536×9=4824" "228×7=1596"
Replace-Text "211×4=844" "777×3=2331"
Replace-Text "195×8=1560" "544×6=3264"
Replace-Text "268×7=1876" "763×6=4578"
Replace-Text "259×8=2072" "964×9=8676"
Replace-Text "682×9=6138" "717×2=1434"
Replace-Text "855×8=6840" "184×5=920"
Replace-Text "642×5=3210" "702×7=4914"
Replace-Text "309×2=618" "587×6=3522"
Replace-Text "923×4=3692" "102×7=714"
Replace-Text "191×3=573" "780×6=4680"
Replace-Text "701×6=4206" "989×5=4945"
Replace-Text "135×2=270" "691×2=1382"
